$d = $word.ActiveDocument

$replacements = @(
    @("58÷7=8, 2", "11÷9=1, 2"),
    @("14÷9=1, 5", "16÷4=4, 0"),
    @("82÷7=11, 5", "38÷5=7, 3"),
    @("20÷3=6, 2", "50÷9=5, 5"),
    @("87÷9=9, 6", "24÷2=12, 0"),
    @("33÷5=6, 3", "17÷7=2, 3"),
    @("14÷2=7, 0", "89÷6=14, 5"),
    @("81÷4=20, 1", "33÷9=3, 6"),
    @("69÷5=13, 4", "54÷2=27, 0"),
    @("30÷9=3, 3", "11÷4=2, 3"),
    @("90÷6=15, 0", "28÷7=4, 0"),
    @("80÷7=11, 3", "25÷6=4, 1"),
    @("86÷2=43, 0", "97÷7=13, 6"),
    @("43÷4=10, 3", "72÷4=18, 0"),
    @("79÷9=8, 7", "14÷7=2, 0"),
    @("71÷9=7, 8", "99÷7=14, 1"),
    @("19÷3=6, 1", "11÷8=1, 3"),
    @("81÷7=11, 4", "16÷8=2, 0"),
    @("71÷5=14, 1", "66÷3=22, 0"),
    @("95÷9=10, 5", "86÷5=17, 1"),
    @("55÷7=7, 6", "59÷2=29, 1"),
    @("77÷7=11, 0", "39÷7=5, 4"),
    @("89÷5=17, 4", "78÷5=15, 3"),
    @("73÷4=18, 1", "96÷7=13, 5"),
    @("57÷2=28, 1", "95÷3=31, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

$d.Save()
